$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: Cash
$ws.Range("B2").Value = 230
$ws.Range("C2").Value = 79.58

# Row 3 becomes "Common Equity" (was row 4), row 4 becomes "Combinations" (was row 3) -- rows swapped
$ws.Range("A3").Value = "Common Equity"
$ws.Range("B3").Value = 22
$ws.Range("C3").Value = 7.61

$ws.Range("A4").Value = "Combinations"
$ws.Range("B4").Value = 17
$ws.Range("C4").Value = 5.88

# Row 5: Cash; Combinations
$ws.Range("B5").Value = 14
$ws.Range("C5").Value = 4.84

# Row 6: Cash; Common Equity
$ws.Range("C6").Value = 1.38

# Row 7: label changes from "Cash; Unknown" to "Cash; Debt"
$ws.Range("A7").Value = "Cash; Debt"
$ws.Range("C7").Value = 0.35

# New row 8: Debt - copy formatting from A7 (style index 1: bordered, centered, bold)
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Debt"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 0.35
